$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsBGD   = $wb.Worksheets.Item("BGDPbES")

# --- 1. BGDPbES sheet: hydro (row 6) guaranteed dispatch now set to a flat 0.35 ---
#     Replace the flowed formulas in B6:AK6 with a literal 0.35 across the whole row.
$wsBGD.Range("B6:AK6").Value = 0.35

# --- 2. About sheet: insert two new note lines (hydro dispatch explanation) ---
#     New text goes in rows 17-18, leaving row 19 blank as a separator, matching
#     the existing note paragraph spacing; the old rows 17-23 shift down to 20-26.
$wsAbout.Rows.Item(17).Insert()
$wsAbout.Rows.Item(17).Insert()
$wsAbout.Rows.Item(17).Insert()

$wsAbout.Range("A17").Value = "For hydro, we guarantee a minimum amount of dispatch and also allow plants to "
$wsAbout.Range("A18").Value = "participate in the energy market for extra capacity"

# --- 3. Update selections / active sheet to match the saved view state ---
$wsBGD.Activate()
$wsBGD.Range("B6:AK6").Select()

$wsAbout.Activate()
$wsAbout.Range("A19").Select()
